$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 3-26 (row 11 unchanged)
$updates = @{
    3  = -2
    4  = -3
    5  = -5
    6  = -1
    7  = 2
    8  = -6
    9  = 2
    10 = -2
    12 = -1
    13 = -2
    14 = -3
    15 = 7
    16 = 8
    17 = -1
    18 = 2
    19 = -2
    20 = 7
    21 = -2
    22 = -2
    23 = -4
    24 = -3
    25 = 1
    26 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
